$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns C:D (the duplicate "Internet Archive" / "Promotion" columns)
$ws.Range("C:D").Delete()

# Update view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D7").Select()

# Adjust tab ratio (window splitter position between sheet tabs and horizontal scrollbar)
$wb.Windows.Item(1).TabRatio = 812
